# Update the Monte-Carlo-derived playoff/record projections
# ("fixed monte carlo for ties").
#
# Two sheets are affected:
#   - "Playoff Odds": columns B..L (place-finish probabilities + "made
#     playoffs" totals) for each team row 2..11.
#   - "Record Odds": column F (playoff chance %), column G (expected final
#     record, now including a tie component like "8.5-5.4-0.1"), and
#     column H (most likely record) for each team row 2..11.

$wb = $excel.ActiveWorkbook

# --- Playoff Odds ------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Playoff Odds")

# Columns B..L, rows 2..11 (one inner array per row, in column order).
$playoffData = @(
    @(76.3, 14.6, 4.7, 3.1, 0.8, 0.3, 0.2, 0, 0, 0, 99.8),
    @(10.5, 28.5, 21.5, 15.5, 10.2, 6.3, 3.9, 2.6, 0.8, 0.2, 92.5),
    @(4.8, 16.9, 18, 17.3, 15.8, 10.6, 7.5, 4.7, 3.1, 1.3, 83.40000000000001),
    @(4.8, 18.3, 17.4, 15.1, 12.3, 11.2, 7.9, 6.6, 3.6, 2.8, 79.09999999999999),
    @(1.7, 9.800000000000001, 16.9, 19.3, 17, 13, 9.4, 6.8, 4.4, 1.7, 77.7),
    @(1.1, 6.6, 10.7, 12.1, 14, 14.7, 14.2, 12.4, 9, 5.2, 59.2),
    @(0.4, 2.9, 4.5, 8, 10.1, 12.1, 15.5, 15.9, 15.6, 15, 38),
    @(0.1, 1, 2.2, 4.2, 8.6, 15.7, 17.1, 17.6, 18, 15.5, 31.8),
    @(0.3, 1.1, 3.5, 4, 6.4, 10.1, 13.8, 17.3, 20.9, 22.6, 25.4),
    @(0, 0.3, 0.6, 1.4, 4.8, 6, 10.5, 16.1, 24.6, 35.7, 13.1)
)

for ($i = 0; $i -lt $playoffData.Length; $i++) {
    $rowVals = $playoffData[$i]
    $rowNum = $i + 2
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws1.Cells.Item($rowNum, $j + 2).Value = $rowVals[$j]
    }
}

# --- Record Odds ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Record Odds")

# Column F (Playoff_Chance_Pct), G (Expected_Final_Record), H (Most_Likely_Record)
$recordOddsData = @(
    @(99.8, "10.8-3.2", "11-3"),
    @(92.5, "8.5-5.4-0.1", "8-6"),
    @(83.39999999999999, "7.9-6.0-0.1", "8-6"),
    @(79.10000000000001, "7.4-6.5", "8-6"),
    @(77.7, "7.2-6.7-0.1", "7-7"),
    @(59.2, "6.4-7.5-0.1", "6-8"),
    @(38, "5.7-8.3-0.1", "6-8"),
    @(31.8, "5.4-8.5-0.1", "6-8"),
    @(25.4, "5.6-8.4-0.1", "6-8"),
    @(13.1, "4.7-9.2-0.1", "5-9")
)

for ($i = 0; $i -lt $recordOddsData.Length; $i++) {
    $rowVals = $recordOddsData[$i]
    $rowNum = $i + 2
    $ws2.Cells.Item($rowNum, 6).Value = $rowVals[0]
    $ws2.Cells.Item($rowNum, 7).Value = $rowVals[1]
    $ws2.Cells.Item($rowNum, 8).Value = $rowVals[2]
}
